$d = $word.ActiveDocument

# 1) Remove the "Meta description" paragraph that currently follows the
#    title ("Play A Night Out Free | Enjoy Wild Bonus Features").
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete() | Out-Null
Write-Host "Removed Meta description paragraph."

# 2) Insert a new bold "Play A Night Out Free | Enjoy Wild Bonus Features"
#    paragraph right before the final (italic) paragraph, matching the
#    document's usual "<empty run/><bold run>" paragraph shape.
$last = $d.Paragraphs.Last
$insertPoint = $d.Range($last.Range.Start, $last.Range.Start)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play A Night Out Free | Enjoy Wild Bonus Features</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($frag) | Out-Null

# InsertXML leaves a blank spacer paragraph behind to force the break;
# remove it now that the split has happened.
$spacer = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$spacer.Range.Delete() | Out-Null
Write-Host "Inserted new bold heading paragraph before the closing paragraph."

# 3) Swap the old italic "Create a feature image…" prompt text for the
#    new meta-description text, keeping the existing italic formatting.
$old = 'Create a feature image for "A Night Out" that features a happy Maya warrior with glasses. The image should be in a cartoon style to complement the retro graphics of the game. The warrior should be in a lively city scene, surrounded by beautiful girls and cocktails, in line with the "vida loca" theme of the game. The image should exude an upbeat, fun vibe, representing the lifestyle that the game embodies. The Maya warrior should be depicted with a big smile on their face, enjoying the night out with friends. Overall, the feature image should be eye-catching, colorful, and representative of the game''s spirit.'
$new = 'In this review, learn why A Night Out is a great free online slot game choice with impressive RTP and fun bonus features. Play now for free.'
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
Write-Host "Replaced closing paragraph text."
